# 13.1.3.xlsx update:
#  - extend the data table with columns for years 2020-2023 (E:H)
#  - fix the "indicator code" text in A1/A4/A5/A6 (1.5.4 -> 13.1.3 relabel
#    already matched the English/Russian headers)
#  - D4 ("484") switches from a text value to a real number, repeated
#    across the new year columns
#  - drop the stale "<selection>" left over from the previous save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels -------------------------------------------------
# A1 kept the wrong "1.5.4 ..." indicator code; replace with the correct
# "13.1.3 ..." Kyrgyz text (B1/C1 already had the right 13.1.3 code and
# are left untouched).
$ws.Range("A1").Value = "13.1.3 Кырсыктардын кооптуулугун азайтуунун улуттук стратегияларына ылайык, кырсыктардын кооптуулугун азайтуунун жергиликтүү стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын үлүшү"

# --- Row 3: add year columns 2020..2023 in E3:H3 (style copied from D3) --
$ws.Range("D3").Copy()
$ws.Range("E3:H3").PasteSpecial(-4122)
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# --- Row 4: "484" becomes a real number, repeated through 2023 -----------
$ws.Range("D4").Copy()
$ws.Range("E4:H4").PasteSpecial(-4122)
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# --- Row 5: proportion (%) series through 2023 ----------------------------
$ws.Range("D5").Copy()
$ws.Range("E5:H5").PasteSpecial(-4122)
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# --- Row 6: count of local governments through 2023 -----------------------
$ws.Range("D6").Copy()
$ws.Range("E6:H6").PasteSpecial(-4122)
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# --- tidy up the leftover selection/clipboard marker from the previous save
$excel.CutCopyMode = $false
$ws.Range("A1").Select()
